$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 291; existing rows 291-306 shift down to 294-309.
$ws.Rows("291:293").Insert()

# Constant column values shared by every record in this sheet/region block.
$marketId = 7
$market   = "Terminal Hortofrutícola Agro Chillán"
$region   = "Ñuble"
$codreg   = 16
$tipo     = "Fruta"
$prodId   = 100103
$prod     = "Frutos de hueso (carozo)"
$catId    = 100103006
$cat      = "Nectarín"
$origen   = "Región de O'Higgins"

# New row 291: June Pearl / Especial
$ws.Cells.Item(291, 1).Value = $marketId
$ws.Cells.Item(291, 2).Value = $market
$ws.Cells.Item(291, 3).Value = $region
$ws.Cells.Item(291, 4).Value = 44610
$ws.Cells.Item(291, 5).Value = $codreg
$ws.Cells.Item(291, 6).Value = $tipo
$ws.Cells.Item(291, 7).Value = $prodId
$ws.Cells.Item(291, 8).Value = $prod
$ws.Cells.Item(291, 9).Value = $catId
$ws.Cells.Item(291, 10).Value = $cat
$ws.Cells.Item(291, 11).Value = "June Pearl"
$ws.Cells.Item(291, 12).Value = "Especial"
$ws.Cells.Item(291, 13).Value = 60
$ws.Cells.Item(291, 14).Value = 14000
$ws.Cells.Item(291, 15).Value = 14000
$ws.Cells.Item(291, 16).Value = 14000
$ws.Cells.Item(291, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(291, 18).Value = $origen
$ws.Cells.Item(291, 19).Value = 933
$ws.Cells.Item(291, 20).Value = 15

# New row 292: June Pearl / Primera
$ws.Cells.Item(292, 1).Value = $marketId
$ws.Cells.Item(292, 2).Value = $market
$ws.Cells.Item(292, 3).Value = $region
$ws.Cells.Item(292, 4).Value = 44610
$ws.Cells.Item(292, 5).Value = $codreg
$ws.Cells.Item(292, 6).Value = $tipo
$ws.Cells.Item(292, 7).Value = $prodId
$ws.Cells.Item(292, 8).Value = $prod
$ws.Cells.Item(292, 9).Value = $catId
$ws.Cells.Item(292, 10).Value = $cat
$ws.Cells.Item(292, 11).Value = "June Pearl"
$ws.Cells.Item(292, 12).Value = "Primera"
$ws.Cells.Item(292, 13).Value = 120
$ws.Cells.Item(292, 14).Value = 12000
$ws.Cells.Item(292, 15).Value = 13000
$ws.Cells.Item(292, 16).Value = 12500
$ws.Cells.Item(292, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(292, 18).Value = $origen
$ws.Cells.Item(292, 19).Value = 833
$ws.Cells.Item(292, 20).Value = 15

# New row 293: June Pearl / Segunda
$ws.Cells.Item(293, 1).Value = $marketId
$ws.Cells.Item(293, 2).Value = $market
$ws.Cells.Item(293, 3).Value = $region
$ws.Cells.Item(293, 4).Value = 44610
$ws.Cells.Item(293, 5).Value = $codreg
$ws.Cells.Item(293, 6).Value = $tipo
$ws.Cells.Item(293, 7).Value = $prodId
$ws.Cells.Item(293, 8).Value = $prod
$ws.Cells.Item(293, 9).Value = $catId
$ws.Cells.Item(293, 10).Value = $cat
$ws.Cells.Item(293, 11).Value = "June Pearl"
$ws.Cells.Item(293, 12).Value = "Segunda"
$ws.Cells.Item(293, 13).Value = 120
$ws.Cells.Item(293, 14).Value = 10000
$ws.Cells.Item(293, 15).Value = 11000
$ws.Cells.Item(293, 16).Value = 10500
$ws.Cells.Item(293, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(293, 18).Value = $origen
$ws.Cells.Item(293, 19).Value = 700
$ws.Cells.Item(293, 20).Value = 15

Write-Host "Rows 291-293 inserted and populated with June Pearl data."
